$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.909.33'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.145.06'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '602.31'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.34'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.47%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.138.56'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.528'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.38'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.466'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.92'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.667.33'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.85%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.979.44'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.144.53'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.86'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '489.57'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.70'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.710'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -2.01%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.29'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.29'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.54%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.44%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.98'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '27.66'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.19%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -6.04%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.50%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.71'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.95'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.49%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0397'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '433.68'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.75%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.31'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.932.87'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.00%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.14%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.42%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.84'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.71%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '120.26'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.08%  '
